$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-12 from 45204 to 45207
$ws.Range("C2:C12").Value = 45207
